$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("shiftsPerWorker")

$ws.Range("B2").Value = 7
$ws.Range("B3").Value = 4
$ws.Range("B5").Value = 5
$ws.Range("B6").Value = 1
$ws.Range("B7").Value = 2
